$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.07838249206543
$ws.Range("B1").Value = 2.29838228225708
$ws.Range("C1").Value = 2.394163608551025
$ws.Range("D1").Value = 3.057729482650757
$ws.Range("E1").Value = 2.803457260131836
